$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C4").Value = -12.02089999999999
$ws.Range("B9").Value = 5.280599999999998
$ws.Range("C9").Value = -13.61009999999999
$ws.Range("D9").Value = -7.948300000000001
$ws.Range("B18").Value = 6.5885
$ws.Range("B20").Value = 9.57499999999999
$ws.Range("C23").Value = -12.3303
$ws.Range("C24").Value = -13.0862
$ws.Range("C26").Value = -12.59460000000001
$ws.Range("B27").Value = 6.591700000000004
$ws.Range("D32").Value = -7.158599999999995
$ws.Range("C34").Value = -11.74530000000001
$ws.Range("C35").Value = -11.78
$ws.Range("D38").Value = -7.958499999999999
$ws.Range("D45").Value = -6.902899999999997
$ws.Range("C48").Value = -11.4339
$ws.Range("D51").Value = -8.556500000000005
$ws.Range("C52").Value = -11.2409
$ws.Range("D57").Value = -8.385499999999999
$ws.Range("D64").Value = -7.223199999999991
$ws.Range("C66").Value = -10.8876
$ws.Range("C67").Value = -10.9641
$ws.Range("B69").Value = 5.315799999999995
$ws.Range("B76").Value = 5.345800000000001
$ws.Range("C80").Value = -13.1729
$ws.Range("B82").Value = 5.870799999999997
$ws.Range("D93").Value = -7.057599999999991
$ws.Range("C99").Value = -12.6866
